# Merge - Opp Test Data, Eng Detail, Add Counterparty - 10 Oct 2025
#
# Applies the data-level changes described by the commit:
#  - CAOUsers: replace the counterparty-list admin contact
#      "Brian Miller" -> "Jennie Stewart"
#  - AddContact: replace a test contact name
#      "Kristian M. Whalen" -> "Alan Test"
#  - NewOpportunityCounterparty: add a new "MassEditComments" column (F)
#      with a sample comment, for the new Add-Counterparty mass edit flow
#  - Restores the various sheet selections / active sheet left behind by
#    the author's last interactive session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# CAOUsers sheet: swap the CAO contact name
# ---------------------------------------------------------------------
$wsCao = $wb.Worksheets.Item("CAOUsers")
$wsCao.Range("A2").Value = "Jennie Stewart"

# ---------------------------------------------------------------------
# NewOpportunityCounterparty sheet: new "MassEditComments" column
# ---------------------------------------------------------------------
$wsCounterparty = $wb.Worksheets.Item("NewOpportunityCounterparty")
$wsCounterparty.Range("F1").Value = "MassEditComments"
$wsCounterparty.Range("F1").Font.Bold = $true
$wsCounterparty.Range("F2").Value = "These are counterparty comments from Counterparty List page"
$wsCounterparty.Columns.Item(6).ColumnWidth = 57.7109375

# ---------------------------------------------------------------------
# AddContact sheet: swap the test contact name used for Contact2
# ---------------------------------------------------------------------
$wsAddContact = $wb.Worksheets.Item("AddContact")
$wsAddContact.Range("F2").Value = "Alan Test"

# ---------------------------------------------------------------------
# Restore the per-sheet selections left over from the author's session.
# The final `Select()` call below (AddContact) also leaves that sheet as
# the active / tabSelected sheet, matching the saved workbook view.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("CAOUsers").Range("J14").Select() | Out-Null
$wb.Worksheets.Item("Opportunities").Range("G23").Select() | Out-Null
$wb.Worksheets.Item("NewOpportunityCounterparty").Range("E28").Select() | Out-Null
$wb.Worksheets.Item("AddContact").Range("F4").Select() | Out-Null
